$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting all existing data right by one column.
$ws.Columns.Item(1).Insert()

# New column A header (row 3 is the header row) -- "Match ID"
$ws.Range("A3").Value = "Match ID"
$ws.Range("A3").Font.Bold = $true

# New column A data values (rows 4-18) -- Match ID = 2
$ws.Range("A4:A18").Value = 2
$ws.Range("A4:A18").Font.Bold = $true

# Totals row (row 19, hidden) also gets the Match ID value, but without bold styling
$ws.Range("A19").Value = 2
$ws.Rows.Item(19).AutoFit()

# Update the active selection to match the saved workbook state
$ws.Range("E23").Select()
